# Add a new bullet point ("A little touch to coding practices") at the end
# of the bullet list on the "Purpose" slide (slide 3), right after
# "Create a GIF animation (Create video highlights)".
#
# TextRange.InsertAfter (rather than overwriting the whole .Text property)
# keeps the existing paragraphs' bullet formatting (buFont/buChar/pPr)
# intact, and the newly-inserted paragraph inherits that same formatting -
# matching native PowerPoint "click at end of list, press Enter, type"
# behaviour.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange
$tr.InsertAfter("`rA little touch to coding practices") | Out-Null
